$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 128, shifting existing rows 128..194 down to 129..195
$ws.Rows.Item(128).Insert()

# Populate the newly inserted row 128 with data
$ws.Range("A128").Value2 = 3
$ws.Range("B128").Value2 = "Femacal de La Calera"
$ws.Range("C128").Value2 = "Coquimbo"
$ws.Range("D128").Value2 = 44455
$ws.Range("E128").Value2 = 5
$ws.Range("F128").Value2 = 100112043
$ws.Range("G128").Value2 = "Pepino ensalada"
$ws.Range("H128").Value2 = "Sin especificar"
$ws.Range("I128").Value2 = "Primera"
$ws.Range("J128").Value2 = 90
$ws.Range("K128").Value2 = 13000
$ws.Range("L128").Value2 = 14000
$ws.Range("M128").Value2 = 13444
$ws.Range("N128").Value2 = "$/caja 70 unidades"
$ws.Range("O128").Value2 = "Región de Arica y Parinacota"
$ws.Range("P128").Value2 = 192
$ws.Range("Q128").Value2 = 70
$ws.Range("R128").Value2 = "Hortaliza"
